$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
# Rows 16-22 hold the "Periodo Mora" (E), "Valor Mora" (F) and
# "Salario Basico" (G) detail for worker 73141410 / ALFONSO MARTIN PEREZ.
# The periods are re-sequenced into ascending order (2203 .. 2209) and the
# mora/salario figures are refreshed with the updated account-statement data.

$data = @(
    @{ Row = 16; Periodo = "2203"; ValorMora = 40000;  SalarioBasico = 1000000 },
    @{ Row = 17; Periodo = "2204"; ValorMora = 40000;  SalarioBasico = 1000000 },
    @{ Row = 18; Periodo = "2205"; ValorMora = 40000;  SalarioBasico = 1000000 },
    @{ Row = 19; Periodo = "2206"; ValorMora = 40000;  SalarioBasico = 1000000 },
    @{ Row = 20; Periodo = "2207"; ValorMora = 40000;  SalarioBasico = 1000000 },
    @{ Row = 21; Periodo = "2208"; ValorMora = 40000;  SalarioBasico = 1000000 },
    @{ Row = 22; Periodo = "2209"; ValorMora = 34666;  SalarioBasico = 1000000 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("E$r").Value = $item.Periodo
    $ws.Range("F$r").Value = $item.ValorMora
    $ws.Range("G$r").Value = $item.SalarioBasico
}
